$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper function to set a text value on a cell while avoiding Excel's
# automatic number/date conversion, and without leaving any residual
# style/number-format changes on the cell (matches original inlineStr cells).
function Set-TextValue($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Row 44 / 45 swap (Maker <-> ApeXProtocol) plus updated price/volume ---
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D44" "3.47"
$ws.Range("E44").Value = "  +19.28%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D45" "2.996.98"
$ws.Range("E45").Value = "  +7.99%  "

# --- Remaining per-row Price (D) / Volume(1h) (E) updates ---
Set-TextValue $ws "D2" "66.769.54"
$ws.Range("E2").Value = "  +8.59%  "
Set-TextValue $ws "D3" "3.484.70"
$ws.Range("E3").Value = "  +12.48%  "
Set-TextValue $ws "D4" "1.00"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws "D5" "188.18"
$ws.Range("E5").Value = "  +13.07%  "
Set-TextValue $ws "D6" "546.32"
$ws.Range("E6").Value = "  +7.53%  "
Set-TextValue $ws "D7" "3.480.17"
$ws.Range("E7").Value = "  +12.37%  "
$ws.Range("E8").Value = "  +4.24%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  +8.45%  "
$ws.Range("E11").Value = "  +19.53%  "
Set-TextValue $ws "D12" "54.43"
$ws.Range("E12").Value = "  +6.46%  "
$ws.Range("E13").Value = "  +10.75%  "
Set-TextValue $ws "D14" "9.30"
$ws.Range("E14").Value = "  +7.23%  "
Set-TextValue $ws "D15" "4.029.12"
$ws.Range("E15").Value = "  +12.34%  "
Set-TextValue $ws "D16" "3.474.31"
$ws.Range("E16").Value = "  +12.39%  "
$ws.Range("E17").Value = "  +7.58%  "
Set-TextValue $ws "D18" "66.737.31"
$ws.Range("E18").Value = "  +9.04%  "
$ws.Range("E19").Value = "  +9.08%  "
$ws.Range("E20").Value = "  +11.16%  "
Set-TextValue $ws "D21" "0.988"
$ws.Range("E21").Value = "  +6.16%  "
Set-TextValue $ws "D22" "421.53"
$ws.Range("E22").Value = "  +19.02%  "
Set-TextValue $ws "D23" "3.89"
$ws.Range("E23").Value = "  +8.25%  "
Set-TextValue $ws "D24" "84.28"
$ws.Range("E24").Value = "  +7.83%  "
$ws.Range("E25").Value = "  +9.76%  "
Set-TextValue $ws "D26" "11.07"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +14.63%  "
Set-TextValue $ws "D28" "6.15"
$ws.Range("E28").Value = "  +0.35%  "
Set-TextValue $ws "D29" "11.85"
$ws.Range("E29").Value = "  +10.31%  "
$ws.Range("E30").Value = "  +12.34%  "
Set-TextValue $ws "D31" "30.03"
$ws.Range("E31").Value = "  +10.44%  "
Set-TextValue $ws "D32" "652.17"
$ws.Range("E32").Value = "  +4.22%  "
Set-TextValue $ws "D33" "6.60"
$ws.Range("E33").Value = "  +6.77%  "
$ws.Range("E34").Value = "  +6.65%  "
Set-TextValue $ws "D35" "0.109"
$ws.Range("E35").Value = "  +9.22%  "
Set-TextValue $ws "D36" "59.29"
$ws.Range("E36").Value = "  +5.94%  "
Set-TextValue $ws "D37" "0.0₃0814"
$ws.Range("E37").Value = "  +23.86%  "
Set-TextValue $ws "D38" "38.28"
$ws.Range("E38").Value = "  +9.94%  "
Set-TextValue $ws "D39" "0.999"
$ws.Range("E39").Value = "  -0.13%  "
Set-TextValue $ws "D40" "0.387"
$ws.Range("E40").Value = "  +6.79%  "
$ws.Range("E41").Value = "  +16.16%  "
Set-TextValue $ws "D42" "3.32"
$ws.Range("E42").Value = "  +16.92%  "
$ws.Range("E43").Value = "  +0.20%  "
Set-TextValue $ws "D46" "2.63"
$ws.Range("E46").Value = "  +6.39%  "
Set-TextValue $ws "D47" "2.87"
$ws.Range("E47").Value = "  +17.61%  "
$ws.Range("E48").Value = "  +10.24%  "
$ws.Range("E49").Value = "  +3.70%  "
Set-TextValue $ws "D50" "8.75"
$ws.Range("E50").Value = "  +19.99%  "
